$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$origStyle = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.850.06"
$ws.Range("D2").Style = $origStyle
$ws.Range("E2").Value = "  -2.46%  "

$origStyle = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.652.85"
$ws.Range("D3").Style = $origStyle
$ws.Range("E3").Value = "  -1.00%  "

$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.18%  "

$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.26"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  -0.57%  "

$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +0.12%  "

$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3890"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -1.89%  "

$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3815"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -2.97%  "

$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "51.52"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -1.58%  "

$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.346"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -3.60%  "

$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.002"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.21%  "

$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08472"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  -1.17%  "

$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "23.97"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -2.50%  "

$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.042"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -3.73%  "

$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.072"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  +1.59%  "

$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001314"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -1.91%  "

$origStyle = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.652.51"
$ws.Range("D17").Style = $origStyle
$ws.Range("E17").Value = "  -0.71%  "

$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "94.12"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("E19").Value = "  -0.31%  "

$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.60"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -5.05%  "

$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.986"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -0.17%  "

$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.002"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.22%  "

$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.68"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -0.68%  "

$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.838.66"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -2.52%  "

$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.431"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -0.98%  "

$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.944"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -4.14%  "

$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.04"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -2.29%  "

$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "154.18"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  -2.23%  "

$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.441"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  +0.11%  "

$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "137.70"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -3.56%  "

$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.803"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  -3.17%  "

$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.491"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -1.89%  "

$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.832.73"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -0.84%  "

$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08173"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -1.39%  "

$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.006"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -5.36%  "

$ws.Range("E36").Value = "  -6.04%  "

$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.606"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -4.57%  "

$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.77"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -4.64%  "

$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2669"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -3.64%  "

$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09125"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -1.62%  "

$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.7557"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -2.30%  "

$ws.Range("E42").Value = "  -2.40%  "

$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.421"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -1.57%  "

$ws.Range("E44").Value = "  -0.80%  "

$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6925"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  -2.89%  "

$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.444"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -4.41%  "

$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.104"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -0.70%  "

$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  +0.09%  "

$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.08279"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -1.80%  "

$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.42"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -2.71%  "

$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.223"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -3.84%  "

